$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose text happens to look numeric must be formatted as Text first,
# otherwise Excel will silently coerce the input into a number.
$textCells = @("F2", "G2", "F3", "G3", "A4", "E4", "F4")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("A2").Value = "sfs"
$ws.Range("B2").Value = "sdfsd"
$ws.Range("C2").Value = "first"
$ws.Range("D2").Value = "asdf"
$ws.Range("E2").Value = " "
$ws.Range("F2").Value = "333"
$ws.Range("G2").Value = "333"

# Row 3
$ws.Range("A3").Value = "f"
$ws.Range("B3").Value = "s"
$ws.Range("C3").Value = "f"
$ws.Range("D3").Value = "f"
$ws.Range("E3").Value = " "
$ws.Range("F3").Value = "3"
$ws.Range("G3").Value = "3"

# Row 4
$ws.Range("A4").Value = "3"
$ws.Range("B4").Value = "f"
$ws.Range("C4").Value = "f"
$ws.Range("D4").Value = "f"
$ws.Range("E4").Value = "3333333333333"
$ws.Range("F4").Value = "333"
$ws.Range("G4").Value = "f"

# Rows 5-7 no longer exist in the data; remove them entirely (shift up).
$ws.Range("A5:G7").EntireRow.Delete()

# Column width tweaks
$ws.Columns.Item(1).ColumnWidth = 16
$ws.Columns.Item(2).ColumnWidth = 12
$ws.Columns.Item(5).ColumnWidth = 15
